$d = $word.ActiveDocument

$pairs = @(
    @("2024-12-04 Wednesday", "2024-12-05 Thursday"),
    @("765÷7=109, 2", "196÷5=39, 1"),
    @("634÷5=126, 4", "923÷3=307, 2"),
    @("686÷7=98, 0", "670÷4=167, 2"),
    @("317÷6=52, 5", "115÷3=38, 1"),
    @("923÷5=184, 3", "325÷9=36, 1"),
    @("543÷5=108, 3", "507÷6=84, 3"),
    @("956÷6=159, 2", "227÷2=113, 1"),
    @("973÷7=139, 0", "887÷5=177, 2"),
    @("577÷2=288, 1", "521÷6=86, 5"),
    @("182÷9=20, 2", "176÷6=29, 2"),
    @("934÷8=116, 6", "169÷3=56, 1"),
    @("919÷8=114, 7", "725÷3=241, 2"),
    @("940÷7=134, 2", "876÷9=97, 3"),
    @("678÷4=169, 2", "814÷7=116, 2"),
    @("931÷3=310, 1", "755÷8=94, 3"),
    @("219÷5=43, 4", "325÷2=162, 1"),
    @("140÷8=17, 4", "611÷5=122, 1"),
    @("694÷3=231, 1", "311÷9=34, 5"),
    @("987÷7=141, 0", "619÷7=88, 3"),
    @("901÷9=100, 1", "613÷3=204, 1"),
    @("421÷2=210, 1", "897÷3=299, 0"),
    @("877÷4=219, 1", "978÷4=244, 2"),
    @("497÷5=99, 2", "418÷2=209, 0"),
    @("933÷6=155, 3", "778÷5=155, 3"),
    @("613÷2=306, 1", "166÷3=55, 1")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
